# Update the YouTube video link shown on the poster slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldUrl = "https://youtu.be/8cISaVsOqs0"
$newUrl = "https://youtu.be/Q-77DMNKo34"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tf = $shape.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text -eq $oldUrl) {
                $tr.Text = $newUrl
            }
        }
    }
}
